# Hindalco price-tracker update (2025-10-20 refresh):
# prepend a new "latest price" row above the current top data row,
# shifting all the existing history down by one row, then fill in
# the newly scraped circular's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row right above the first data row (row 2, below
# the frozen header) and copy row 3's (the old row 2's) formatting
# into it so the new row looks like the rest of the data rows instead
# of inheriting the bold header style.
$ws.Rows("2:2").Insert()
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# Fill in the newly published circular's values.
$ws.Range("A2").Value = 44
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 279.75
$ws.Range("E2").Value = "18.10.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-18-october-2025.pdf"

# The row that used to be "F24" (12.08.2025 circular) now sits at F25
# after the insert above, and it picks up a hyperlink for the first
# time, matching the PDF link text already present in that cell.
$ws.Hyperlinks.Add($ws.Range("F25"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf")
# Adding the hyperlink applies Excel's built-in blue/underlined
# "Hyperlink" style; restore the plain centered style used by every
# other data cell in that column (copy format from the row below).
$ws.Range("F26").Copy()
$ws.Range("F25").PasteSpecial(-4122)
